# Updated results, model selection, main effects and multiple comparisons
#
# Rewrites the GLMM main-effects table: two new model terms (length,
# day_night) were added, expanding the 3-term model (habitat_type, season,
# habitat_type:season) into a full 4-way factorial main-effects +
# interactions table (4 main effects, 6 two-way, 4 three-way, 1 four-way
# = 15 rows), and refreshes every statistic/df/p.value accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row is unchanged, but set it explicitly for robustness.
$ws.Cells.Item(1,1).Value = "term"
$ws.Cells.Item(1,2).Value = "statistic"
$ws.Cells.Item(1,3).Value = "df"
$ws.Cells.Item(1,4).Value = "p.value"

$ws.Cells.Item(2,1).Value = "habitat_type"
$ws.Cells.Item(2,2).Value = 84.78774276261
$ws.Cells.Item(2,3).Value = 4
$ws.Cells.Item(2,4).Value = 0.0000000000000000168268952037243

$ws.Cells.Item(3,1).Value = "season"
$ws.Cells.Item(3,2).Value = 31.791892300358
$ws.Cells.Item(3,3).Value = 3
$ws.Cells.Item(3,4).Value = 0.000000578954057181146

$ws.Cells.Item(4,1).Value = "length"
$ws.Cells.Item(4,2).Value = 2.17781887118537
$ws.Cells.Item(4,3).Value = 1
$ws.Cells.Item(4,4).Value = 0.140012764067774

$ws.Cells.Item(5,1).Value = "day_night"
$ws.Cells.Item(5,2).Value = 65.4392079502853
$ws.Cells.Item(5,3).Value = 3
$ws.Cells.Item(5,4).Value = 0.0000000000000404023285817303

$ws.Cells.Item(6,1).Value = "habitat_type:season"
$ws.Cells.Item(6,2).Value = 419.041298153316
$ws.Cells.Item(6,3).Value = 12
$ws.Cells.Item(6,4).Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000349727223196716

$ws.Cells.Item(7,1).Value = "habitat_type:length"
$ws.Cells.Item(7,2).Value = 8.46945745134344
$ws.Cells.Item(7,3).Value = 4
$ws.Cells.Item(7,4).Value = 0.075818443192257

$ws.Cells.Item(8,1).Value = "season:length"
$ws.Cells.Item(8,2).Value = 3.02202919600533
$ws.Cells.Item(8,3).Value = 3
$ws.Cells.Item(8,4).Value = 0.388241162466469

$ws.Cells.Item(9,1).Value = "habitat_type:day_night"
$ws.Cells.Item(9,2).Value = 76.1396917712486
$ws.Cells.Item(9,3).Value = 12
$ws.Cells.Item(9,4).Value = 0.0000000000223610556311327

$ws.Cells.Item(10,1).Value = "season:day_night"
$ws.Cells.Item(10,2).Value = 103.689354843691
$ws.Cells.Item(10,3).Value = 9
$ws.Cells.Item(10,4).Value = 0.00000000000000000281640857090181

$ws.Cells.Item(11,1).Value = "length:day_night"
$ws.Cells.Item(11,2).Value = 0.882059051648772
$ws.Cells.Item(11,3).Value = 3
$ws.Cells.Item(11,4).Value = 0.829755071728027

$ws.Cells.Item(12,1).Value = "habitat_type:season:length"
$ws.Cells.Item(12,2).Value = 73.9655600013696
$ws.Cells.Item(12,3).Value = 10
$ws.Cells.Item(12,4).Value = 0.00000000000756111626058715

$ws.Cells.Item(13,1).Value = "habitat_type:season:day_night"
$ws.Cells.Item(13,2).Value = 92.7968813686353
$ws.Cells.Item(13,3).Value = 35
$ws.Cells.Item(13,4).Value = 0.000000395979420761007

$ws.Cells.Item(14,1).Value = "habitat_type:length:day_night"
$ws.Cells.Item(14,2).Value = 24.6280047359372
$ws.Cells.Item(14,3).Value = 12
$ws.Cells.Item(14,4).Value = 0.0166874110907106

$ws.Cells.Item(15,1).Value = "season:length:day_night"
$ws.Cells.Item(15,2).Value = 31.6771227555644
$ws.Cells.Item(15,3).Value = 9
$ws.Cells.Item(15,4).Value = 0.000226400760352013

$ws.Cells.Item(16,1).Value = "habitat_type:season:length:day_night"
$ws.Cells.Item(16,2).Value = 43.854964813563
$ws.Cells.Item(16,3).Value = 29
$ws.Cells.Item(16,4).Value = 0.0378799099043248
